# Add a new order row (row 8) to the order list sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text values - safe to assign directly, Excel will reuse/create
# shared-string entries as needed.
$ws.Range("A8").Value = "D-100"
$ws.Range("B8").Value = "NTU"
$ws.Range("C8").Value = "Fries, Chicken Nugget"
$ws.Range("D8").Value = "Fries : spicy"
$ws.Range("F8").Value = "NEW"
$ws.Range("G8").Value = "Cash"

# The "takeaway" column value is the literal text "false". Assigning that
# string directly would be auto-coerced to a boolean by Excel, so instead
# copy it from an existing cell that already holds the same text value,
# which preserves it as a shared string.
$ws.Range("E5").Copy()
$ws.Range("E8").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Newly written cells inherit the worksheet's default column style; strip
# that back to the workbook's base "Normal" style so the new row matches
# the unstyled rows above it (rows 5-7) and styles.xml stays unchanged.
$ws.Range("A8:G8").Style = "Normal"
